{"js": "// Change the heading text in the first table cell from\n// \"In the County Court Money Claims Centre \" to \"In the County Court\".\n// The original paragraph is split across three runs: \"In the Coun\" + \"t\" +\n// \"y Court Money Claims Centre \" (the last run carries a trailing space).\n// We only need to edit the text of that third run, turning it into\n// \"y Court\" (search on the full run text, including the trailing space,\n// so nothing is left behind).\n\nconst searchResults = context.document.body.search(\"y Court Money Claims Centre \", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (const range of searchResults.items) {\n  range.insertText(\"y Court\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Change the heading text in the first table cell from\n# \"In the County Court Money Claims Centre \" to \"In the County Court\"\n# by deleting the trailing \" Money Claims Centre \" portion.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \" Money Claims Centre \"\n$find.Replacement.Text = \"\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue\n\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)  # wdReplaceAll\n\n$d.Save()\n"}
